$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cliente")

# Row 9: replace client name and CPF (CPF becomes a plain number, not text)
$ws.Range("B9").Value = "Aloisia Andrade"
$ws.Range("D9").Value = 44444444444

# Remove the duplicate trailing row (row 11 - "Jailson Silva" repeat)
$ws.Rows.Item(11).Delete()

# Update the active selection to reflect the edited area
$ws.Range("D10").Select()
